$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace Acceptance Criteria text (column I, rows 3-27) with shortened versions ---
$ws.Range("I3").Value2 = "Los usuarios pueden registrarse, iniciar sesión y tener roles definidos."
$ws.Range("I4").Value2 = "Cada cliente puede registrar sus datos y los de sus mascotas."
$ws.Range("I5").Value2 = "Roles creados, editables, asignables a cuentas."
$ws.Range("I6").Value2 = "Formulario de mascota funcional con validaciones obligatorias."
$ws.Range("I7").Value2 = "Contraseña actualizable, último acceso registrado correctamente."
$ws.Range("I8").Value2 = "Solo registros activos pueden iniciar sesión o ser asignados."
$ws.Range("I9").Value2 = "Direcciones con latitud/longitud y opción “principal”."
$ws.Range("I10").Value2 = "Historial visible por fechas y editable por el nutricionista."
$ws.Range("I11").Value2 = "Los productos y platos pueden crearse, combinarse y publicarse."
$ws.Range("I12").Value2 = "Permite elegir productos, cantidades y marcar “es crudo”."
$ws.Range("I13").Value2 = "Lista de platos publicados con precio, imagen y descripción."
$ws.Range("I14").Value2 = "Carga en menos de 2 segundos promedio."
$ws.Range("I15").Value2 = "Permitir agendar consultas y generar dietas personalizadas."
$ws.Range("I16").Value2 = "Selección de fecha, mascota y nutricionista disponibles."
$ws.Range("I17").Value2 = "Dieta con platos, fechas, instrucciones y frecuencia."
$ws.Range("I18").Value2 = "Encriptación y validación de permisos en consultas."
$ws.Range("I19").Value2 = "El cliente puede realizar pedidos y confirmar entregas."
$ws.Range("I20").Value2 = "Pedido genera detalle y total automáticamente."
$ws.Range("I21").Value2 = "Confirmación con marca “entregado” y fecha registrada."
$ws.Range("I22").Value2 = "Subtotal correcto según cantidad y precio."
$ws.Range("I23").Value2 = "Registrar pagos y sus historiales."
$ws.Range("I24").Value2 = "Pasarela seleccionable, monto y referencia válidos."
$ws.Range("I25").Value2 = "Pagos listados por fecha, estado y método."
$ws.Range("I26").Value2 = "Promedio < 2 segundos."
$ws.Range("I27").Value2 = "Hash seguro y política de contraseñas."

# --- Drop the wrapText style (cellXf 8) from I3:I27, reverting to default formatting ---
$ws.Range("I3:I27").Style = "Normal"

# --- Column width adjustments ---
$ws.Columns("C").ColumnWidth = 64.5
$ws.Columns("I").ColumnWidth = 90.6

# --- Add 5 new blank rows (28-32) with the same row height as the data rows ---
$ws.Rows("28:32").RowHeight = 40.05
$ws.Range("H28").Style = "Normal"
$ws.Range("H29").Style = "Normal"

# --- View: zoom + scroll position + selection ---
$excel.ActiveWindow.Zoom = 74
$ws.Range("H28").Select()
